# submit vote from tags in comments
# Add a new "Handle" column (E) that derives an @-handle from the
# employee name in column B, for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("E1").Value = "Handle"

# First data row gets its own (non-shared) formula, matching how Excel
# records the first cell of a fill/autofill operation.
$ws.Range("E2").Formula = '="@"&B2'

# The remaining rows (E3:E11) are filled with the same relative formula;
# assigning the formula to the whole range at once makes the host engine
# record it as a shared formula group, same as Excel's fill-down does.
$ws.Range("E3:E11").Formula = '="@"&B3'

# Leave the selection on the newly filled range, mirroring the end state
# of a fill-down action that finishes on the last row.
$ws.Range("E2:E11").Select()
